$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as of the new scrape.
# Numeric-looking text values are forced to remain text (matching the
# original inline-string cell contents) by temporarily applying a text
# number format, then resetting the cell style back to Normal so no
# stray formatting is left behind.

$ws.Range("D2").Value = '45.224.41'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '2.364.42'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.73%  '
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.612'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.61%  '
$ws.Range("E11").Value = '  -2.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.48'
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.980'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.33%  '
$ws.Range("D15").Value = '2.725.25'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.53%  '
$ws.Range("D17").Value = '2.356.49'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").Value = '45.186.96'
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.27'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +14.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.22%  '
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.86%  '
$ws.Range("E25").Value = '  -1.33%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.58%  '
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.34%  '
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0946'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '36.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '168.83'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.19%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.116'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.69'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.90'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.97'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0354'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("D44").Value = '1.878.77'
$ws.Range("E44").Value = '  +14.35%  '
$ws.Range("E45").Value = '  -5.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.20%  '
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.93'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.33'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.06%  '
